# Apply schedule updates: resize day-columns and move/add/remove the
# inline "assignment" strings on sheets s1 / s2 / Therapists / Patients
# to reflect the new "taller_verano" workshop session.

$wb = $excel.ActiveWorkbook

# Excel's Range.ColumnWidth is expressed in "characters" and Excel adds a
# constant ~5/6 character padding when it writes the <col width=".."/>
# attribute in the OOXML. Subtract that offset so the saved width comes
# out to the exact integer we want.
$colPad = 5 / 6
function Set-StoredColWidth($ws, $colIndex, $storedWidth) {
    $ws.Columns.Item($colIndex).ColumnWidth = $storedWidth - $colPad
}

$KINE   = "kine_1 | agu | kine:javi | s1 | n=1"
$FONO   = "fono_1 | agu | fono:maca | s2 | n=1"
$TALLER = "taller_verano | agu | fono:maca, kine:javi | s2 | n=1"

# ---------------------------------------------------------------------
# Sheet 1: "s1"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("s1")

Set-StoredColWidth $ws1 2 37
Set-StoredColWidth $ws1 5 12

$ws1.Range("D5").ClearContents()
$ws1.Range("C5").Value = $KINE

$ws1.Range("C10").ClearContents()
$ws1.Range("E10").ClearContents()
$ws1.Range("B10").Value = $KINE
$ws1.Range("D10").Value = $KINE

# ---------------------------------------------------------------------
# Sheet 2: "s2"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("s2")

Set-StoredColWidth $ws2 2 55
Set-StoredColWidth $ws2 3 55
Set-StoredColWidth $ws2 4 55
Set-StoredColWidth $ws2 5 12

$ws2.Range("B3").Value = $TALLER

$ws2.Range("D5").Value = $TALLER

$ws2.Range("D6").ClearContents()

$ws2.Range("B7").Value = $FONO

$ws2.Range("E8").ClearContents()
$ws2.Range("D8").Value = $FONO

$ws2.Range("C10").Value = $TALLER

# ---------------------------------------------------------------------
# Sheet 3: "Therapists"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Therapists")

Set-StoredColWidth $ws3 3 55
Set-StoredColWidth $ws3 4 55

$ws3.Range("C3").Value = $TALLER
$ws3.Range("D3").Value = $TALLER

$ws3.Range("D7").Value = $FONO

$ws3.Range("C10").Value = $KINE

$ws3.Range("C14").Value = $KINE

$ws3.Range("C19").Value = $TALLER
$ws3.Range("D19").Value = $TALLER

$ws3.Range("C23").Value = $TALLER
$ws3.Range("D23").Value = $TALLER

$ws3.Range("D24").ClearContents()

$ws3.Range("D26").Value = $FONO

$ws3.Range("C28").Value = $KINE

$ws3.Range("D35").ClearContents()

$ws3.Range("C37").ClearContents()

# ---------------------------------------------------------------------
# Sheet 4: "Patients"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Patients")

Set-StoredColWidth $ws4 3 55

$ws4.Range("C3").Value = $TALLER

$ws4.Range("C7").Value = $FONO

$ws4.Range("C10").Value = $KINE

$ws4.Range("C14").Value = $KINE

$ws4.Range("C19").Value = $TALLER

$ws4.Range("C23").Value = $TALLER

$ws4.Range("C24").ClearContents()

$ws4.Range("C26").Value = $FONO

$ws4.Range("C28").Value = $KINE

$ws4.Range("C35").ClearContents()

$ws4.Range("C37").ClearContents()
